$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# The "0-非、1-是" remark on the three Y/N flag fields (IsRelated/IsLnrelNear/IsLimit,
# rows 11-13, column G) is replaced with a two-line "0:非" / "1:是" note.
$note = "0:非" + [char]10 + "1:是"
$ws.Range("G11").Value = $note
$ws.Range("G12").Value = $note
$ws.Range("G13").Value = $note

# The cell style already wraps text, so the now two-line content makes Excel grow
# those rows to 32.4pt (two lines at the sheet's normal 16.2pt row height).
$ws.Rows("11:13").RowHeight = 32.4

# Leave the cursor on G14, matching the saved selection in the workbook.
$ws.Range("G14").Select()
